$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (deterministic) - only Execution Time changes
$ws.Range("H2").Value = 0.00300288200378418

# Row 3 (randomized)
$ws.Range("C3").Value = 515.53
$ws.Range("D3").Value = 514.91
$ws.Range("E3").Value = 0.62
$ws.Range("F3").Value = 1030.44
$ws.Range("G3").Value = 515.22
$ws.Range("H3").Value = 1.626104116439819

# Row 4 (deterministic) - only Execution Time changes
$ws.Range("H4").Value = 0.002505064010620117

# Row 5 (randomized)
$ws.Range("C5").Value = 648.26
$ws.Range("D5").Value = 480.69
$ws.Range("E5").Value = 167.56
$ws.Range("F5").Value = 1128.95
$ws.Range("G5").Value = 564.47
$ws.Range("H5").Value = 1.693790912628174

# Row 6 (deterministic) - only Execution Time changes
$ws.Range("H6").Value = 0.006505250930786133

# Row 7 (randomized)
$ws.Range("C7").Value = 555.7
$ws.Range("D7").Value = 548.53
$ws.Range("E7").Value = 7.17
$ws.Range("F7").Value = 1653.52
$ws.Range("G7").Value = 551.17
$ws.Range("H7").Value = 3.580146789550781

# Row 8 (deterministic) - only Execution Time changes
$ws.Range("H8").Value = 0.007514238357543945

# Row 9 (randomized)
$ws.Range("C9").Value = 644.21
$ws.Range("D9").Value = 508.83
$ws.Range("E9").Value = 135.37
$ws.Range("F9").Value = 1796.22
$ws.Range("G9").Value = 598.74
$ws.Range("H9").Value = 3.922915458679199

# Row 10 (deterministic) - only Execution Time changes
$ws.Range("H10").Value = 0.01050400733947754

# Row 11 (randomized)
$ws.Range("C11").Value = 541.97
$ws.Range("D11").Value = 523.75
$ws.Range("E11").Value = 18.22
$ws.Range("F11").Value = 2141.73
$ws.Range("G11").Value = 535.4299999999999
$ws.Range("H11").Value = 5.896767854690552

# Row 12 (deterministic) - only Execution Time changes
$ws.Range("H12").Value = 0.01150727272033691

# Row 13 (randomized)
$ws.Range("C13").Value = 696.2
$ws.Range("D13").Value = 467.39
$ws.Range("E13").Value = 228.81
$ws.Range("F13").Value = 2328.68
$ws.Range("G13").Value = 582.17
$ws.Range("H13").Value = 6.262850284576416
